# "made proxy for individual heat PJ"
#
# This script:
#  1. Adds a new comment on E22 (proxy note explaining the data issue).
#  2. Fills in row 22 (A22:H22) with a new data row:
#       remind | SSP2-Base | Wind | DK | consumption|individual heat|individual heat use | PJ | 53.7 | 70.1
#     (this also introduces a new shared string for the "variables" column).
#  3. Copies the number-format/font styling used by the neighbouring E19 cell
#     (0.0000 number format, Arial font) onto E22:E24, matching the style
#     used for other "variables" cells in this unit column.
#  4. Moves the sheet's active selection to H23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Comment on E22 -----------------------------------------------------
$ws.Range("E22").AddComment("Data in source doesn’t add up, and unit is unclear. Must be wrong?") | Out-Null

# --- 3. Style E22:E24 like E19 (numFmt 0.0000 / Arial) before writing values
# so the style index gets reused/created before we touch the cell values.
$ws.Range("E19").Copy() | Out-Null
$ws.Range("E22:E24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 2. Row 22 values -------------------------------------------------------
$ws.Range("A22").Value = "remind"
$ws.Range("B22").Value = "SSP2-Base"
$ws.Range("C22").Value = "Wind"
$ws.Range("D22").Value = "DK"
$ws.Range("E22").Value = "consumption|individual heat|individual heat use"
$ws.Range("F22").Value = "PJ"
$ws.Range("G22").Value = 53.7
$ws.Range("H22").Value = 70.1

# --- 4. Update active selection --------------------------------------------
$ws.Range("H23").Select() | Out-Null
